# Update cryptocurrency price/volume/hour data (refresh run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. Values must remain TEXT
# (matching the source data, which stores numbers/percents as strings),
# so each cell is forced to text format before the write and restored
# to the default "Normal" style afterward so no stray styling is left behind.
$updates = [ordered]@{
    'D2' = '314.45'
    'E2' = '0.52%'
    'G2' = '21'
    'D3' = '37.38'
    'E3' = '-0.74%'
    'G3' = '21'
    'D4' = '5.141'
    'E4' = '-0.36%'
    'G4' = '21'
    'D5' = '0.07933'
    'E5' = '0.58%'
    'G5' = '21'
    'D6' = '8.449'
    'E6' = '1.85%'
    'G6' = '21'
    'D7' = '1.896'
    'E7' = '-1.03%'
    'G7' = '21'
    'E8' = '4.35%'
    'G8' = '21'
    'D9' = '0.9327'
    'E9' = '1.22%'
    'G9' = '21'
    'D10' = '0.1273'
    'E10' = '4.79%'
    'G10' = '21'
    'D11' = '0.1925'
    'E11' = '-0.09%'
    'G11' = '21'
    'D12' = '0.08939'
    'E12' = '-2.64%'
    'G12' = '21'
    'D13' = '0.03387'
    'E13' = '1.48%'
    'G13' = '21'
    'D14' = '0.09515'
    'E14' = '-1.00%'
    'G14' = '21'
    'D15' = '0.001371'
    'E15' = '-0.56%'
    'G15' = '21'
    'D16' = '0.006174'
    'E16' = '6.95%'
    'G16' = '21'
    'D17' = '3.394'
    'E17' = '-3.46%'
    'G17' = '21'
    'D18' = '4.444'
    'E18' = '0.70%'
    'G18' = '21'
    'E19' = '1.38%'
    'G19' = '21'
    'D20' = '6.486'
    'E20' = '23.19%'
    'G20' = '21'
    'D21' = '0.1302'
    'E21' = '2.29%'
    'G21' = '21'
    'D22' = '0.2303'
    'E22' = '-11.13%'
    'G22' = '21'
    'D23' = '0.04353'
    'E23' = '-0.39%'
    'G23' = '21'
    'D24' = '0.001199'
    'E24' = '-3.97%'
    'G24' = '21'
    'D25' = '0.004227'
    'E25' = '-10.30%'
    'G25' = '21'
    'D26' = '0.0001325'
    'E26' = '8.39%'
    'G26' = '21'
    'D27' = '0.0003965'
    'E27' = '-98.11%'
    'G27' = '21'
    'G28' = '21'
    'G29' = '21'
    'G30' = '21'
    'G31' = '21'
    'G32' = '21'
    'G33' = '21'
    'G34' = '21'
    'G35' = '21'
    'G36' = '21'
    'G37' = '21'
    'G38' = '21'
    'D39' = '0.02290'
    'E39' = '-0.91%'
    'G39' = '21'
    'D40' = '0.05122'
    'E40' = '0.60%'
    'G40' = '21'
    'D41' = '0.007472'
    'E41' = '0.24%'
    'G41' = '21'
    'E42' = '2.19%'
    'G42' = '21'
    'D43' = '0.008547'
    'E43' = '-5.51%'
    'G43' = '21'
    'D44' = '0.002062'
    'E44' = '5.60%'
    'G44' = '21'
    'D45' = '0.007937'
    'E45' = '-7.84%'
    'G45' = '21'
    'D46' = '0.00006326'
    'E46' = '-4.58%'
    'G46' = '21'
    'E47' = '-0.43%'
    'G47' = '21'
    'D48' = '0.002855'
    'G48' = '21'
    'D49' = '0.001683'
    'E49' = '40.13%'
    'G49' = '21'
    'D50' = '0.00002092'
    'E50' = '-0.43%'
    'G50' = '21'
    'E51' = '-0.43%'
    'G51' = '21'
}

foreach ($addr in $updates.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$addr]
    $rng.Style = "Normal"
}

